$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds the "Förändrad" (last-changed) date for each data row.
# Determine the last used row in column C (xlUp = -4162) starting from the bottom.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 367 }

# Update every data row (row 2 through the last used row) in column C
# from 45205 (2023-10-06) to 45206 (2023-10-07).
$ws.Range("C2:C$lastRow").Value = 45206
